# Updated SRS Review workbook to V1.9
#  - Row 9  : mark "Accepted", add comment about moving the Context diagram
#  - Row 10 : mark "Accepted", add comment about reverting Req_..._013
#  - Row 11 : mark "Accepted", add comment about splitting Req_..._001 / 018
#  - Move the active selection to H11 / scroll the view up a couple of rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainRun($rng) {
    # Re-assert the "normal" (non-bold) run formatting so the emitted run
    # gets an explicit <rPr>, matching the rest of the sheet's rich-text runs.
    $rng.Font.Bold = $false
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 11
}

# ---------------------------------------------------------------------------
# Row 9 - "Moved the Context diagram under the project description"
# ---------------------------------------------------------------------------
$e9 = $ws.Range("E9")
$e9.VerticalAlignment = -4108   # xlVAlignCenter
$e9.WrapText = $true
$e9.Value = "Accepted"

$ws.Range("G9").Value = "Moved the Context diagram under the project description"

$ws.Rows(9).RowHeight = 28.8

# ---------------------------------------------------------------------------
# Row 10 - "Reverted back Req_PO1_DGC_SRS_013_V01 to its original state..."
# ---------------------------------------------------------------------------
$e10 = $ws.Range("E10")
$e10.VerticalAlignment = -4108
$e10.WrapText = $true
$e10.Value = "Accepted"

$text10 = "Reverted back Req_PO1_DGC_SRS_013_V01 to its original state, which idicates that the software design should follow the Flow Chart"
$g10 = $ws.Range("G10")
$g10.Value = $text10
$g10.Characters(15, 23).Font.Bold = $true
Set-PlainRun ($g10.Characters(38, 92))

# row 10's height (43.2) is unchanged

# ---------------------------------------------------------------------------
# Row 11 - "Made Req_PO1_DGC_SRS_001_V01 more specific ... Req_..._018 ..."
# ---------------------------------------------------------------------------
$e11 = $ws.Range("E11")
$e11.VerticalAlignment = -4108
$e11.WrapText = $true
$e11.Value = "Accepted"

$text11 = "Made Req_PO1_DGC_SRS_001_V01 more specific and only concerned with operands 1 and 2, and the operator.`nReq_PO1_DGC_SRS_018_V01 is already concerened with storing/calculatung the result."
$g11 = $ws.Range("G11")
$g11.Value = $text11
$g11.Characters(6, 23).Font.Bold = $true
Set-PlainRun ($g11.Characters(29, 75))
$g11.Characters(104, 23).Font.Bold = $true
Set-PlainRun ($g11.Characters(127, 59))

$ws.Rows(11).RowHeight = 86.4

# ---------------------------------------------------------------------------
# Extend the "Accepted/Rejected" conditional formatting + data validation
# from E2:E10 down to E2:E11 now that row 11 also has a disposition.
# ---------------------------------------------------------------------------
$e10fcs = $ws.Range("E10").FormatConditions
for ($i = 1; $i -le $e10fcs.Count; $i++) {
    $e10fcs.Item($i).ModifyAppliesToRange($ws.Range("E10:E11")) | Out-Null
}

$ws.Range("E2:E11").Validation.Delete()
$ws.Range("E2:E11").Validation.Add(3, 1, 1, '"Accepted, Rejected"') | Out-Null

# ---------------------------------------------------------------------------
# Selection / scroll position
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H11").Select() | Out-Null

Write-Output "edit complete"
